$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("E","G","H","I","J","K","M","N","O","P","Q","R","S","T")

$data = @{
    2 = @(3,11.27436066666667,33.823082,0.6286406327466463,0.6286406327466463,3,24.59985066666667,73.79955200000001,0.4261128628549158,0.4261128628549157,277.3475887621404,2496.128298859264,0.2678718597265992,0.2678718597265991)
    3 = @(3,11.27436066666667,33.823082,0.6286406327466463,0.6286406327466463,3,2.577566,7.732698000000001,0.04464799572729766,0.04464799572729766,29.06040872613734,261.543678535236,0.02806754428487796,0.02806754428487796)
    4 = @(3,11.27436066666667,33.823082,0.6286406327466463,0.6286406327466463,3,0.4574283333333333,1.372285,0.007923466662817385,0.007923466662817385,5.157212009152222,46.41490808237,0.004981013096460479,0.004981013096460479)
    5 = @(3,11.27436066666667,33.823082,0.6286406327466463,0.6286406327466463,3,30.09598833333334,90.28796500000001,0.5213156747549692,0.5213156747549692,339.3130270897923,3053.817243808131,0.3277202156387087,0.3277202156387087)
    6 = @(3,2.395357333333333,7.186071999999999,0.1335613605242407,0.1335613605242408,3,24.59985066666667,73.79955200000001,0.4261128628549158,0.4261128628549157,58.92543269330489,530.328894239744,0.05691221369978175,0.05691221369978176)
    7 = @(3,2.395357333333333,7.186071999999999,0.1335613605242407,0.1335613605242408,3,2.577566,7.732698000000001,0.04464799572729766,0.04464799572729766,6.174191620250667,55.56772458225601,0.005963247054018362,0.005963247054018364)
    8 = @(3,2.395357333333333,7.186071999999999,0.1335613605242407,0.1335613605242408,3,0.4574283333333333,1.372285,0.007923466662817385,0.007923466662817385,1.095704312724444,9.86133881452,0.001058268987554355,0.001058268987554356)
    9 = @(3,2.395357333333333,7.186071999999999,0.1335613605242407,0.1335613605242408,3,30.09598833333334,90.28796500000001,0.5213156747549692,0.5213156747549692,72.09064635816445,648.8158172234801,0.06962763078288627,0.06962763078288628)
    10 = @(3,1.512471333333333,4.537414,0.08433302464847796,0.08433302464847797,3,24.59985066666667,73.79955200000001,0.4261128628549158,0.4261128628549157,37.20656893761422,334.859120438528,0.03593538656617712,0.03593538656617712)
    11 = @(3,1.512471333333333,4.537414,0.08433302464847796,0.08433302464847797,3,2.577566,7.732698000000001,0.04464799572729766,0.04464799572729766,3.898494684774667,35.086452162972,0.003765300524175332,0.003765300524175333)
    12 = @(3,1.512471333333333,4.537414,0.08433302464847796,0.08433302464847797,3,0.4574283333333333,1.372285,0.007923466662817385,0.007923466662817385,0.6918472412211111,6.22662517099,0.0006682099093767719,0.000668209909376772)
    13 = @(3,1.512471333333333,4.537414,0.08433302464847796,0.08433302464847797,3,30.09598833333334,90.28796500000001,0.5213156747549692,0.5213156747549692,45.51931960250111,409.6738764225101,0.04396412764874873,0.04396412764874874)
    14 = @(3,2.752319,8.256957,0.1534649820806351,0.1534649820806351,3,24.59985066666667,73.79955200000001,0.4261128628549158,0.4261128628549157,67.70663638702933,609.359727483264,0.06539340286235776,0.06539340286235776)
    15 = @(3,2.752319,8.256957,0.1534649820806351,0.1534649820806351,3,2.577566,7.732698000000001,0.04464799572729766,0.04464799572729766,7.094283875554001,63.848554879986,0.006851903864226006,0.006851903864226007)
    16 = @(3,2.752319,8.256957,0.1534649820806351,0.1534649820806351,3,0.4574283333333333,1.372285,0.007923466662817385,0.007923466662817385,1.258988692971667,11.330898236745,0.001215974669425779,0.001215974669425779)
    17 = @(3,2.752319,8.256957,0.1534649820806351,0.1534649820806351,3,30.09598833333334,90.28796500000001,0.5213156747549692,0.5213156747549692,82.83376051361168,745.5038446225051,0.08000370068462552,0.08000370068462553)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cellRef = "$($cols[$i])$row"
        $ws.Range($cellRef).Value = $vals[$i]
    }
}

Write-Host "Done updating cells"
